$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Phase 1" planning rows added beneath the existing To-Do list
$ws.Range("B12").Value = "Phase 1:"
$ws.Range("B13").Value = "Modify Player Mov Speed"
$ws.Range("B14").Value = "Change From fire to energy Fields"
$ws.Range("B15").Value = "Change Behavior of the orbs"
$ws.Range("B16").Value = "Change The points system"
$ws.Range("C16").Value = "Player will have to fill the score bar to advance to phase 2"
$ws.Range("B18").Value = "Define Phase 2"

# Columns B and C need to widen to fit the new (longer) text
$ws.Columns("B").ColumnWidth = 27.65
$ws.Columns("C").ColumnWidth = 47.8

# Move the selection down to the newly added last entry
$ws.Range("B18").Select()
